# Applies the "added attendance and fees page" edit described by the diff:
#   - renames the sole sheet "A" -> "Sheet1" (and bumps its internal sheetId,
#     matching what real Excel does when you recreate/duplicate a sheet)
#   - appends a new student record (row 2) to the roster table
#   - bumps the shared base font from 11pt to 12pt
#   - grows the header/new row to the resulting 15.6pt row height
#   - records the new used range + suppresses the "number stored as text"
#     advisory for the age/grade columns, which are intentionally text

$wb = $excel.ActiveWorkbook

# --- Rename "A" -> "Sheet1", while giving the sheet a fresh sheetId -------
# A plain .Name= rename keeps the original sheetId. The source workbook's
# sheetId moved from 1 to 2, which is what happens when the sheet is
# duplicated and the original is deleted (the duplicate gets a new id).
# Copy-before-self preserves all existing content/formatting of sheet "A".
$wsOld = $wb.Worksheets.Item(1)
$wsOld.Copy($wsOld)

# The copy is now the first sheet; the original "A" got pushed to slot 2.
$wb.Worksheets.Item(2).Delete()
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"

# --- Bump the workbook's base font size 11 -> 12 --------------------------
$wb.Styles.Item("Normal").Font.Size = 12

# --- Append the new student record (row 2) ---------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Krish"
$ws.Range("C2").Value = "vu4f2324074@pvppcoe.ac.in"

# Age/grade are stored as text ("18"/"10"), not numbers, in the source data.
# Force text storage, then drop back to the Normal style so no stray
# per-cell number format lingers on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "18"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "aids"

# --- Row heights follow the larger 12pt font -------------------------------
$ws.Rows.Item(1).RowHeight = 15.6
$ws.Rows.Item(2).RowHeight = 15.6

# --- Suppress the "number stored as text" advisory over the used range ----
$ws.Range("A1:F2").Errors.Item(1).Ignore = $true
